# ShoppingTrip.xlsx -- add named ranges + a new "Sheet1" with a number
# grid and a couple of SUM/AVERAGE formulas that lean on the new names.

$wb = $excel.ActiveWorkbook
$shoppingTrip = $wb.Worksheets.Item("Shopping Trip")

# ---------------------------------------------------------------------
# 1. New worksheet "Sheet1", placed after "Shopping Trip"
# ---------------------------------------------------------------------
$sheet1 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $shoppingTrip)
$sheet1.Name = "Sheet1"

# D1:M15 number grid -- row r, column offset j (0-based from D) => r + 2*j
for ($r = 1; $r -le 15; $r++) {
    for ($j = 0; $j -le 9; $j++) {
        $col = 4 + $j
        $sheet1.Cells.Item($r, $col).Value = $r + (2 * $j)
    }
}

# ---------------------------------------------------------------------
# 2. Defined names
# ---------------------------------------------------------------------
$wb.Names.Add("Taxes", "='Shopping Trip'!`$C`$2:`$C`$6")
$wb.Names.Add("Taxplus", "='Shopping Trip'!`$C`$2:`$D`$6")
$wb.Names.Add("Imanuel", "='Shopping Trip'!`$B`$4,'Shopping Trip'!`$D`$6,'Shopping Trip'!`$D`$2")
$wb.Names.Add("Mike", "=Sheet1!`$D`$1:`$M`$15")
$wb.Names.Add("Zack", "=Sheet1!`$D`$1:`$M`$15")

# ---------------------------------------------------------------------
# 3. New formulas
# ---------------------------------------------------------------------
$shoppingTrip.Range("G8").Formula = "=SUM(Taxplus)"
$shoppingTrip.Range("G10").Formula = "=SUM(Imanuel)"

$sheet1.Range("R10").Formula = "=SUM(Zack)"
$sheet1.Range("R11").Formula = "=AVERAGE(Zack)"

# ---------------------------------------------------------------------
# 4. Selections / active sheet to mirror the saved view state
# ---------------------------------------------------------------------
[void]$shoppingTrip.Range("B2:B6").Select()
[void]$sheet1.Range("Q16").Select()
[void]$sheet1.Activate()
